# Apply updated evaluation metric values across the three worksheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3558718861209965
$wsSummary.Range("C2").Value = 0.06510416666666667
$wsSummary.Range("D2").Value = 0.8928571428571429
$wsSummary.Range("E2").Value = 0.1213592233009709
$wsSummary.Range("F2").Value = 0.2520161290322581
$wsSummary.Range("G2").Value = 0.5996309963099631
$wsSummary.Range("H2").Value = 0.8013643659711074
$wsSummary.Range("I2").Value = 25
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 3

# --- Sheet 2: Classification Report ---
$wsClassification = $wb.Worksheets.Item("Classification Report")

# Row 2 ("0")
$wsClassification.Range("B2").Value = 0.9831460674157303
$wsClassification.Range("C2").Value = 0.3277153558052435
$wsClassification.Range("D2").Value = 0.4915730337078651

# Row 3 ("1")
$wsClassification.Range("B3").Value = 0.06510416666666667
$wsClassification.Range("C3").Value = 0.8928571428571429
$wsClassification.Range("D3").Value = 0.1213592233009709

# Row 4 ("accuracy")
$wsClassification.Range("B4").Value = 0.3558718861209965
$wsClassification.Range("C4").Value = 0.3558718861209965
$wsClassification.Range("D4").Value = 0.3558718861209965
$wsClassification.Range("E4").Value = 0.3558718861209965

# Row 5 ("macro avg")
$wsClassification.Range("B5").Value = 0.5241251170411985
$wsClassification.Range("C5").Value = 0.6102862493311931
$wsClassification.Range("D5").Value = 0.306466128504418

# Row 6 ("weighted avg")
$wsClassification.Range("B6").Value = 0.937407325029656
$wsClassification.Range("C6").Value = 0.3558718861209965
$wsClassification.Range("D6").Value = 0.4731282175310092

# --- Sheet 3: Confusion Matrix ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 ("Actual 0")
$wsConfusion.Range("B2").Value = 175
$wsConfusion.Range("C2").Value = 359

# Row 3 ("Actual 1")
$wsConfusion.Range("B3").Value = 3
$wsConfusion.Range("C3").Value = 25
